$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (LinearRegression) - values change slightly
$ws.Range("B2").Value = -24051.46218472848
$ws.Range("C2").Value = -24051.46218472848
$ws.Range("D2").Value = -24051.46218472848

# Row 3 (RandomForestRegressor) - values become negative
$ws.Range("B3").Value = -103.7432752122634
$ws.Range("C3").Value = -148.8889595641804
$ws.Range("D3").Value = -2452.143655514816

# Row 4 - name changes from GradientBoostingRegressor to DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.7912295017441876
$ws.Range("C4").Value = 0.7941126771330138
$ws.Range("D4").Value = -8976.015052853105

# Row 5 - name changes from AdaBoostRegressor to MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = -67.00292919653369
$ws.Range("C5").Value = -101.2861620935063
$ws.Range("D5").Value = -10793.11343175781
